$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force plain-text storage for numeric-looking price values (matches source data as text)
foreach ($addr in @("D4","D5","D6","D7","D8","D9","D10","D11","D13","D14","D15","D16","D17","D18","D19","D20","D22","D23","D25","D26","D27","D28","D29","D30","D31","D32","D33","D34","D35","D36","D38","D39","D40","D41","D43","D44","D45","D46","D47","D48","D49","D50","D51")) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply cell value updates from the diff
$ws.Range("D2").Value = "26.341.65"
$ws.Range("E2").Value = "  -3.85%  "
$ws.Range("D3").Value = "1.767.59"
$ws.Range("E3").Value = "  -2.74%  "
$ws.Range("D4").Value = "0.9966"
$ws.Range("E4").Value = "  -0.80%  "
$ws.Range("B5").Value = "USDC"
$ws.Range("C5").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D5").Value = "0.9956"
$ws.Range("E5").Value = "  -0.70%  "
$ws.Range("B6").Value = "BNB"
$ws.Range("C6").Value = "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
$ws.Range("D6").Value = "304.84"
$ws.Range("E6").Value = "  -2.33%  "
$ws.Range("D7").Value = "0.4269"
$ws.Range("E7").Value = "  +1.07%  "
$ws.Range("D8").Value = "0.3634"
$ws.Range("E8").Value = "  +1.44%  "
$ws.Range("D9").Value = "0.07159"
$ws.Range("E9").Value = "  -0.08%  "
$ws.Range("D10").Value = "0.8461"
$ws.Range("E10").Value = "  -0.45%  "
$ws.Range("D11").Value = "20.32"
$ws.Range("E11").Value = "  +0.14%  "
$ws.Range("D12").Value = "1.750.69"
$ws.Range("E12").Value = "  -5.27%  "
$ws.Range("D13").Value = "5.234"
$ws.Range("E13").Value = "  -1.90%  "
$ws.Range("D14").Value = "6.416"
$ws.Range("E14").Value = "  +0.40%  "
$ws.Range("D15").Value = "0.06865"
$ws.Range("E15").Value = "  -0.88%  "
$ws.Range("D16").Value = "1.001"
$ws.Range("E16").Value = "  -0.60%  "
$ws.Range("D17").Value = "78.72"
$ws.Range("E17").Value = "  -3.39%  "
$ws.Range("D18").Value = "0.000008661"
$ws.Range("E18").Value = "  -2.12%  "
$ws.Range("D19").Value = "0.9956"
$ws.Range("E19").Value = "  -0.61%  "
$ws.Range("D20").Value = "14.99"
$ws.Range("E20").Value = "  -1.60%  "
$ws.Range("D21").Value = "26.333.00"
$ws.Range("E21").Value = "  -4.72%  "
$ws.Range("D22").Value = "5.088"
$ws.Range("E22").Value = "  +0.03%  "
$ws.Range("D23").Value = "11.07"
$ws.Range("E23").Value = "  +0.73%  "
$ws.Range("D24").Value = "1.978.71"
$ws.Range("E24").Value = "  -4.46%  "
$ws.Range("D25").Value = "151.46"
$ws.Range("E25").Value = "  -1.55%  "
$ws.Range("D26").Value = "1.852"
$ws.Range("E26").Value = "  -6.27%  "
$ws.Range("D27").Value = "18.01"
$ws.Range("E27").Value = "  -1.30%  "
$ws.Range("D28").Value = "5.076"
$ws.Range("E28").Value = "  -0.68%  "
$ws.Range("D29").Value = "113.64"
$ws.Range("E29").Value = "  +0.15%  "
$ws.Range("D30").Value = "1.800"
$ws.Range("E30").Value = "  +3.03%  "
$ws.Range("D31").Value = "0.08923"
$ws.Range("E31").Value = "  +0.39%  "
$ws.Range("D32").Value = "0.7273"
$ws.Range("E32").Value = "  -1.96%  "
$ws.Range("D33").Value = "1.120"
$ws.Range("E33").Value = "  +0.72%  "
$ws.Range("D34").Value = "4.313"
$ws.Range("E34").Value = "  -3.82%  "
$ws.Range("D35").Value = "0.9948"
$ws.Range("E35").Value = "  -0.86%  "
$ws.Range("D36").Value = "2.723"
$ws.Range("E36").Value = "  -8.33%  "
$ws.Range("E37").Value = "  +1.24%  "
$ws.Range("D38").Value = "0.05141"
$ws.Range("E38").Value = "  -1.29%  "
$ws.Range("D39").Value = "0.01878"
$ws.Range("E39").Value = "  -1.64%  "
$ws.Range("D40").Value = "0.4913"
$ws.Range("E40").Value = "  -2.11%  "
$ws.Range("D41").Value = "0.1608"
$ws.Range("E41").Value = "  -2.03%  "
$ws.Range("D43").Value = "6.291"
$ws.Range("E43").Value = "  +0.06%  "
$ws.Range("D44").Value = "8.006"
$ws.Range("E44").Value = "  -3.05%  "
$ws.Range("D45").Value = "104.68"
$ws.Range("E45").Value = "  -0.83%  "
$ws.Range("D46").Value = "10.15"
$ws.Range("E46").Value = "  -2.66%  "
$ws.Range("D47").Value = "0.9942"
$ws.Range("E47").Value = "  -0.66%  "
$ws.Range("B48").Value = "Cronos"
$ws.Range("C48").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D48").Value = "0.06186"
$ws.Range("E48").Value = "  -3.89%  "
$ws.Range("B49").Value = "Decentraland"
$ws.Range("C49").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D49").Value = "0.4482"
$ws.Range("E49").Value = "  -3.07%  "
$ws.Range("D50").Value = "1.614"
$ws.Range("E50").Value = "  +0.68%  "
$ws.Range("D51").Value = "1.737"
$ws.Range("E51").Value = "  +3.42%  "
